# Update numeric values in column F for the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1138
    $ws.Range("F6").Value = 137
    $ws.Range("F10").Value = 5142
    $ws.Range("F11").Value = 4760
}
